$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

$ws.Range("A111").Value = "Estado De México"

$ws.Range("B17").Value = "Amatenango De La Frontera"
$ws.Range("B18").Value = "Amatenango Del Valle"
$ws.Range("B21").Value = "Bejucal De Ocampo"
$ws.Range("B25").Value = "Chiapa De Corzo"
$ws.Range("A84").Value = "Ciudad De México"
$ws.Range("B88").Value = "Cuajimalpa De Morelos"
$ws.Range("B108").Value = "San Juan Del Río"
$ws.Range("B111").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B113").Value = "Almoloya De Alquisiras"
$ws.Range("B114").Value = "Almoloya De Juárez"
$ws.Range("B115").Value = "Almoloya Del Río"
$ws.Range("B120").Value = "Atizapán De Zaragoza"
$ws.Range("B126").Value = "Chapa De Mota"
$ws.Range("B128").Value = "Coacalco De Berriozábal"
$ws.Range("B134").Value = "Ecatepec De Morelos"
$ws.Range("B138").Value = "Ixtapan De La Sal"
$ws.Range("B150").Value = "Naucalpan De Juárez"
$ws.Range("B157").Value = "San Antonio La Isla"
$ws.Range("B158").Value = "San Felipe Del Progreso"
$ws.Range("B159").Value = "San Simón De Guerrero"
$ws.Range("B161").Value = "Soyaniquilpan De Juárez"
$ws.Range("B169").Value = "Tenango Del Valle"
$ws.Range("B178").Value = "Tlalnepantla De Baz"
$ws.Range("B184").Value = "Valle De Bravo"
$ws.Range("B185").Value = "Valle De Chalco Solidaridad"
$ws.Range("B186").Value = "Villa De Allende"
$ws.Range("B187").Value = "Villa Del Carbón"
$ws.Range("B199").Value = "San Miguel De Allende"
$ws.Range("B200").Value = "Apaseo El Alto"
$ws.Range("B201").Value = "Apaseo El Grande"
$ws.Range("B207").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B211").Value = "Jaral Del Progreso"
$ws.Range("B218").Value = "Purísima Del Rincón"
$ws.Range("B221").Value = "San Diego De La Unión"
$ws.Range("B223").Value = "San Francisco Del Rincón"
$ws.Range("B224").Value = "San Luis De La Paz"
$ws.Range("B225").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B227").Value = "Silao De La Victoria"
$ws.Range("B231").Value = "Valle De Santiago"
$ws.Range("B236").Value = "Acapulco De Juárez"
$ws.Range("B239").Value = "Ajuchitlán Del Progreso"
$ws.Range("B240").Value = "Alcozauca De Guerrero"
$ws.Range("B243").Value = "Atenango Del Río"
$ws.Range("B244").Value = "Atlamajalcingo Del Monte"
$ws.Range("B245").Value = "Atoyac De Álvarez"
$ws.Range("B246").Value = "Ayutla De Los Libres"
$ws.Range("B248").Value = "Buenavista De Cuéllar"
$ws.Range("B249").Value = "Chilapa De Álvarez"
$ws.Range("B250").Value = "Chilpancingo De Los Bravo"
$ws.Range("B255").Value = "Coyuca De Benítez"
$ws.Range("B256").Value = "Coyuca De Catalán"
$ws.Range("B260").Value = "Cuetzala Del Progreso"
$ws.Range("B261").Value = "Cutzamala De Pinzón"
$ws.Range("B267").Value = "Huitzuco De Los Figueroa"
$ws.Range("B268").Value = "Iguala De La Independencia"
$ws.Range("B269").Value = "Ixcateopan De Cuauhtémoc"
$ws.Range("B272").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B285").Value = "Taxco De Alarcón"
$ws.Range("B287").Value = "Técpan De Galeana"
$ws.Range("B289").Value = "Tepecoacuilco De Trujano"
$ws.Range("B291").Value = "Tixtla De Guerrero"
$ws.Range("B294").Value = "Tlalixtaquilla De Maldonado"
$ws.Range("B295").Value = "Tlapa De Comonfort"
$ws.Range("B307").Value = "Agua Blanca De Iturbide"
$ws.Range("B311").Value = "Atotonilco El Grande"
$ws.Range("B316").Value = "Huasca De Ocampo"
$ws.Range("B319").Value = "Jacala De Ledezma"
$ws.Range("B322").Value = "Mineral Del Chico"
$ws.Range("B323").Value = "Mineral Del Monte"
$ws.Range("B324").Value = "Mixquiahuala De Juárez"
$ws.Range("B325").Value = "Molango De Escamilla"
$ws.Range("B327").Value = "Nopala De Villagrán"
$ws.Range("B328").Value = "Omitlán De Juárez"
$ws.Range("B329").Value = "Pachuca De Soto"
$ws.Range("B331").Value = "Progreso De Obregón"
$ws.Range("B337").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B340").Value = "Tenango De Doria"
$ws.Range("B342").Value = "Tepehuacán De Guerrero"
$ws.Range("B343").Value = "Tezontepec De Aldama"
$ws.Range("B347").Value = "Tula De Allende"
$ws.Range("B348").Value = "Tulancingo De Bravo"
$ws.Range("B349").Value = "Zacualtipán De Ángeles"
$ws.Range("B350").Value = "Zapotlán De Juárez"
$ws.Range("B353").Value = "Ahualulco De Mercado"
$ws.Range("B356").Value = "Atotonilco El Alto"
$ws.Range("B357").Value = "Autlán De Navarro"
$ws.Range("B361").Value = "Encarnación De Díaz"
$ws.Range("B363").Value = "Huejuquilla El Alto"
$ws.Range("B364").Value = "Ixtlahuacán De Los Membrillos"
$ws.Range("B368").Value = "Lagos De Moreno"
$ws.Range("B371").Value = "Ojuelos De Jalisco"
$ws.Range("B373").Value = "San Cristóbal De La Barranca"
$ws.Range("B374").Value = "San Martín De Bolaños"
$ws.Range("B377").Value = "Tamazula De Gordiano"
$ws.Range("B380").Value = "Tepatitlán De Morelos"
$ws.Range("B382").Value = "Tizapán El Alto"
$ws.Range("B388").Value = "Zapotlán El Grande"
$ws.Range("B404").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B466").Value = "Coatlán Del Río"
$ws.Range("B475").Value = "Puente De Ixtla"
$ws.Range("B479").Value = "Tetela Del Volcán"
$ws.Range("B480").Value = "Tlaltizapán De Zapata"
$ws.Range("B485").Value = "Zacualpan De Amilpas"
$ws.Range("B488").Value = "Ixtlán Del Río"
$ws.Range("B506").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B508").Value = "Ayoquezco De Aldama"
$ws.Range("B513").Value = "Coicoyán De Las Flores"
$ws.Range("B514").Value = "El Barrio De La Soledad"
$ws.Range("B515").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B516").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B517").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B518").Value = "Huautla De Jiménez"
$ws.Range("B519").Value = "Ixtlán De Juárez"
$ws.Range("B520").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B523").Value = "Mariscala De Juárez"
$ws.Range("B527").Value = "Oaxaca De Juárez"
$ws.Range("B528").Value = "Ocotlán De Morelos"
$ws.Range("B529").Value = "Putla Villa De Guerrero"
$ws.Range("B540").Value = "San Antonino El Alto"
$ws.Range("B546").Value = "San Dionisio Del Mar"
$ws.Range("B553").Value = "San José Del Progreso"
$ws.Range("B556").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B634").Value = "Tanetze De Zaragoza"
$ws.Range("B636").Value = "Tlacolula De Matamoros"
$ws.Range("B638").Value = "Villa De Chilapa De Díaz"
$ws.Range("B639").Value = "Villa De Etla"
$ws.Range("B640").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B641").Value = "Villa De Zaachila"
$ws.Range("B661").Value = "Ayotoxco De Guerrero"
$ws.Range("B672").Value = "Chila De La Sal"
$ws.Range("B677").Value = "Cuapiaxtla De Madero"
$ws.Range("B681").Value = "Cuayuca De Andrade"
$ws.Range("B682").Value = "Cuetzalan Del Progreso"
$ws.Range("B694").Value = "Huehuetlán El Chico"
$ws.Range("B698").Value = "Huitzilan De Serdán"
$ws.Range("B700").Value = "Ixcamilpa De Guerrero"
$ws.Range("B703").Value = "Izúcar De Matamoros"
$ws.Range("B710").Value = "Los Reyes De Juárez"
$ws.Range("B717").Value = "Palmar De Bravo"
$ws.Range("B735").Value = "San Nicolás De Los Ranchos"
$ws.Range("B738").Value = "San Salvador El Seco"
$ws.Range("B739").Value = "San Salvador El Verde"
$ws.Range("B744").Value = "Tecali De Herrera"
$ws.Range("B750").Value = "Tepanco De López"
$ws.Range("B751").Value = "Tepango De Rodríguez"
$ws.Range("B752").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B757").Value = "Tepexi De Rodríguez"
$ws.Range("B758").Value = "Tetela De Ocampo"
$ws.Range("B763").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B774").Value = "Tuzamapan De Galeana"
$ws.Range("B788").Value = "Amealco De Bonfil"
$ws.Range("B789").Value = "Cadereyta De Montes"
$ws.Range("B793").Value = "Jalpan De Serra"
$ws.Range("B794").Value = "Pinal De Amoles"
$ws.Range("B797").Value = "San Juan Del Río"
$ws.Range("B807").Value = "Ciudad Del Maíz"
$ws.Range("B811").Value = "Mexquitic De Carmona"
$ws.Range("B814").Value = "Santa María Del Río"
$ws.Range("B815").Value = "Soledad De Graciano Sánchez"
$ws.Range("B817").Value = "Villa De Arista"
$ws.Range("B818").Value = "Villa De Guadalupe"
$ws.Range("B819").Value = "Villa De Reyes"
$ws.Range("B855").Value = "Soto La Marina"
$ws.Range("B859").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B863").Value = "Contla De Juan Cuamatzi"
$ws.Range("B868").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B869").Value = "Nanacamilpa De Mariano Arista"
$ws.Range("B874").Value = "San Pablo Del Monte"
$ws.Range("B877").Value = "Tepetitla De Lardizábal"
$ws.Range("B880").Value = "Tetla De La Solidaridad"
$ws.Range("B893").Value = "Amatlán De Los Reyes"
$ws.Range("B896").Value = "Camarón De Tejeda"
$ws.Range("B900").Value = "Castillo De Teayo"
$ws.Range("B909").Value = "Cosamaloapan De Carpio"
$ws.Range("B919").Value = "Hueyapan De Ocampo"
$ws.Range("B920").Value = "Ignacio De La Llave"
$ws.Range("B922").Value = "Ixhuatlán De Madero"
$ws.Range("B923").Value = "Ixhuatlán Del Café"
$ws.Range("B929").Value = "Juchique De Ferrer"
$ws.Range("B932").Value = "Lerdo De Tejada"
$ws.Range("B933").Value = "Martínez De La Torre"
$ws.Range("B934").Value = "Medellín De Bravo"
$ws.Range("B937").Value = "Mixtla De Altamirano"
$ws.Range("B944").Value = "Paso De Ovejas"
$ws.Range("B945").Value = "Paso Del Macho"
$ws.Range("B948").Value = "Poza Rica De Hidalgo"
$ws.Range("B952").Value = "Sayula De Alemán"
$ws.Range("B953").Value = "Soledad De Doblado"
$ws.Range("B967").Value = "Vega De Alatorre"
$ws.Range("B972").Value = "Zozocolco De Hidalgo"
$ws.Range("B983").Value = "Teúl De González Ortega"
$ws.Range("B984").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B985").Value = "Villa De Cos"

$ws.Range("D143").Value = 0.009503372164316373
$ws.Range("D303").Value = 0.0994788473329246
$ws.Range("D772").Value = 0.009503372164316373

$ws.Rows("992:996").Delete()
